$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.451.32"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "2.684.51"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("D5").Value = "'523.64"
$ws.Range("D6").Value = "'144.65"
$ws.Range("E6").Value = "  -1.32%  "
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").Value = "'0.574"
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("D9").Value = "2.708.01"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").Value = "'6.41"
$ws.Range("E10").Value = "  +2.15%  "
$ws.Range("E11").Value = "  -1.85%  "
$ws.Range("D12").Value = "'0.338"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("E13").Value = "  +1.94%  "
$ws.Range("D14").Value = "3.156.82"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").Value = "60.477.61"
$ws.Range("E15").Value = "  +1.80%  "
$ws.Range("D16").Value = "'21.29"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "2.697.89"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").Value = "'349.78"
$ws.Range("E19").Value = "  -2.77%  "
$ws.Range("D20").Value = "'4.52"
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("D22").Value = "'6.29"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "'63.10"
$ws.Range("E24").Value = "  +1.47%  "
$ws.Range("D25").Value = "'0.420"
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("D26").Value = "'0.170"
$ws.Range("E26").Value = "  +4.80%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").Value = "0.0₃0817"
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("D29").Value = "'7.32"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").Value = "'6.89"
$ws.Range("E30").Value = "  +7.50%  "
$ws.Range("D31").Value = "'0.998"
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("D32").Value = "'19.17"
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").Value = "'147.97"
$ws.Range("E34").Value = "  -1.86%  "
$ws.Range("D35").Value = "'4.25"
$ws.Range("E35").Value = "  +4.10%  "
$ws.Range("D36").Value = "'0.962"
$ws.Range("E36").Value = "  -3.67%  "
$ws.Range("E37").Value = "  +7.42%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "'0.873"
$ws.Range("E38").Value = "  +1.48%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'1.52"
$ws.Range("E39").Value = "  +6.34%  "
$ws.Range("D40").Value = "'36.92"
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").Value = "'3.65"
$ws.Range("E41").Value = "  -2.47%  "
$ws.Range("D42").Value = "'283.17"
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'0.997"
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.611"
$ws.Range("E45").Value = "  -2.08%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.0987"
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("D47").Value = "2.136.41"
$ws.Range("E47").Value = "  +5.97%  "
$ws.Range("D48").Value = "'4.91"
$ws.Range("E48").Value = "  +2.60%  "
$ws.Range("D49").Value = "'0.0538"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("D50").Value = "'0.0234"
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("D51").Value = "'10.45"
$ws.Range("E51").Value = "  +1.62%  "
